$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The "Comentarios" row is the last row (row 11). We insert a new row
# just above it for the "Radioisótopos / Afectados" entry.
$comentariosRow = $t.Rows.Item($t.Rows.Count)
$newRow = $t.Rows.Add($comentariosRow)

# --- Column 1 (label cell): two paragraphs "Radioisótopos" / "Afectados"
$cell1 = $newRow.Cells.Item(1)
$cell1.Range.Text = "Radioisótopos" + [char]13 + "Afectados"

# --- Column 2 (value cell): bracketed instruction text
$cell2 = $newRow.Cells.Item(2)
$text2 = "[Indicar cuáles de los siguientes radioisótopos podrían ser liberados a la atmósfera (Uruanio 238/Torio 232/Radio 226/Polonio 218)]"
$cell2.Range.Text = $text2

# The paragraph mark of the value cell keeps an underline formatting
# remnant while the run text itself stays without underline.
$p2 = $t.Cell($t.Rows.Count - 1, 2).Range.Paragraphs.Item(1)
$p2.Range.Font.Underline = 1
$valueTextRange = $d.Range($p2.Range.Start, $p2.Range.Start + $text2.Length)
$valueTextRange.Font.Underline = 0

# Word leaves a "_GoBack" bookmark marking the last edited spot.
$labelCell = $t.Cell($t.Rows.Count - 1, 1)
$goBackRange = $d.Range($labelCell.Range.Start, $labelCell.Range.Start + "Radioisótopos".Length)
$d.Bookmarks.Add("_GoBack", $goBackRange)
